$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("N1_D40")
$ws1.Range("E2").Value = 0.017
$ws1.Range("F2").Value = 10.48
$ws1.Range("E3").Value = 0.016
$ws1.Range("F3").Value = 10.54
$ws1.Range("E4").Value = 0.016
$ws1.Range("F4").Value = 10.54
$ws1.Range("E5").Value = 0.016
$ws1.Range("F5").Value = 10.51
$ws1.Range("E6").Value = 0.017
$ws1.Range("F6").Value = 10.54
$ws1.Range("E7").Value = 0.016
$ws1.Range("F7").Value = 10.46
$ws1.Range("E8").Value = 0.017
$ws1.Range("E9").Value = 0.017
$ws1.Range("F9").Value = 10.48
$ws1.Range("E10").Value = 0.017
$ws1.Range("F10").Value = 10.54
$ws1.Range("E11").Value = 0.016
$ws1.Range("F11").Value = 10.48
$ws1.Range("E12").Value = 0.0165
$ws1.Range("F12").Value = 10.505

$ws2 = $wb.Worksheets.Item("N1_D60")
$ws2.Range("E2").Value = 0.025
$ws2.Range("F2").Value = 9.99
$ws2.Range("E3").Value = 0.025
$ws2.Range("F3").Value = 9.890000000000001
$ws2.Range("E4").Value = 0.025
$ws2.Range("F4").Value = 10.04
$ws2.Range("E5").Value = 0.025
$ws2.Range("F5").Value = 10.02
$ws2.Range("E6").Value = 0.025
$ws2.Range("F6").Value = 9.890000000000001
$ws2.Range("E7").Value = 0.025
$ws2.Range("F7").Value = 9.94
$ws2.Range("E8").Value = 0.025
$ws2.Range("F8").Value = 10
$ws2.Range("E9").Value = 0.025
$ws2.Range("F9").Value = 9.99
$ws2.Range("E10").Value = 0.027
$ws2.Range("F10").Value = 10.02
$ws2.Range("E11").Value = 0.025
$ws2.Range("F11").Value = 10.03
$ws2.Range("E12").Value = 0.0252
$ws2.Range("F12").Value = 9.980999999999998

$ws3 = $wb.Worksheets.Item("N1_D80")
$ws3.Range("E2").Value = 0.038
$ws3.Range("F2").Value = 13.51
$ws3.Range("E3").Value = 0.038
$ws3.Range("F3").Value = 13.55
$ws3.Range("E4").Value = 0.039
$ws3.Range("F4").Value = 13.51
$ws3.Range("E5").Value = 0.038
$ws3.Range("F5").Value = 13.5
$ws3.Range("E6").Value = 0.039
$ws3.Range("F6").Value = 13.57
$ws3.Range("E7").Value = 0.038
$ws3.Range("F7").Value = 13.59
$ws3.Range("F8").Value = 13.56
$ws3.Range("E9").Value = 0.039
$ws3.Range("F9").Value = 13.65
$ws3.Range("E10").Value = 0.039
$ws3.Range("F10").Value = 13.57
$ws3.Range("E11").Value = 0.039
$ws3.Range("F11").Value = 13.56
$ws3.Range("E12").Value = 0.0387
$ws3.Range("F12").Value = 13.557

$ws4 = $wb.Worksheets.Item("N1_D100")
$ws4.Range("E2").Value = 0.057
$ws4.Range("F2").Value = 11.61
$ws4.Range("E3").Value = 0.056
$ws4.Range("F3").Value = 11.61
$ws4.Range("E4").Value = 0.056
$ws4.Range("E5").Value = 0.056
$ws4.Range("F5").Value = 11.69
$ws4.Range("E6").Value = 0.056
$ws4.Range("F6").Value = 11.61
$ws4.Range("E7").Value = 0.056
$ws4.Range("F7").Value = 11.55
$ws4.Range("E8").Value = 0.057
$ws4.Range("F8").Value = 11.63
$ws4.Range("E9").Value = 0.056
$ws4.Range("F9").Value = 11.6
$ws4.Range("E10").Value = 0.057
$ws4.Range("F10").Value = 11.55
$ws4.Range("E11").Value = 0.057
$ws4.Range("F11").Value = 11.62
$ws4.Range("E12").Value = 0.05640000000000001
$ws4.Range("F12").Value = 11.608
